# daily auto push: 2026-02-26 10:03 UTC
#
# A new reading was appended to the "sei1" log and needs to be inserted in
# its chronological slot: a new row for 2026/02/26 17:00 (rank 201) belongs
# right before the 2026/12/29 block, i.e. at row 889. Inserting it there
# pushes every following row (old 889-930) down by one, so the sheet grows
# from A1:D930 to A1:D931.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Push rows 889..930 down one slot to make room for the new entry.
$ws.Rows(889).Insert()

# Copy the date/weekday text from the row right above (same date,
# 2026/02/26 is a 木=Thursday) instead of re-typing the literal string, so
# Excel stores it as plain text instead of auto-converting it to a date
# serial number (matching how the rest of the column is stored).
$ws.Range("A888:B888").Copy($ws.Range("A889"))

# Fill in this entry's own time-of-day and ranking values.
$ws.Range("C889").Value = 17
$ws.Range("D889").Value = 201

Write-Host "Inserted 2026/02/26 17:00 (rank 201) at row 889; rows 889-930 shifted to 890-931."
